# Reorders rows 2-13 of the data table: the columns that vary per-record
# (Fecha, Volumen, Precio mínimo/máximo/promedio, Origen, Precio $/Kg) are
# shuffled into a new row order, while the columns that are constant for
# every record (Mercado ID, Mercado, Región, Codreg, Tipo, Producto...,
# Unidad de comercialización, Kg / unidad) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source (old) row number, for rows 2..13.
$rowMap = @{
    2  = 13
    3  = 5
    4  = 3
    5  = 8
    6  = 4
    7  = 7
    8  = 12
    9  = 11
    10 = 10
    11 = 9
    12 = 6
    13 = 2
}

# Capture the "before" values of the columns that move (D, M, N, O, P, R, S)
# for every source row first, since several rows are both a source and a
# destination and we must not read an already-overwritten value.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    $ws.Cells.Item($destRow, 4).Value2 = $vals.D
    $ws.Cells.Item($destRow, 13).Value2 = $vals.M
    $ws.Cells.Item($destRow, 14).Value2 = $vals.N
    $ws.Cells.Item($destRow, 15).Value2 = $vals.O
    $ws.Cells.Item($destRow, 16).Value2 = $vals.P
    $ws.Cells.Item($destRow, 18).Value2 = $vals.R
    $ws.Cells.Item($destRow, 19).Value2 = $vals.S
}
